# Add 2022-Q3 data
# 1) Insert a new row of summary data ("2022-Q3") at the top of the "总计"
#    (totals) sheet's data table, pushing the existing quarters down by one
#    row.
# 2) Insert a brand-new worksheet named "2022-Q3" right after "总计",
#    populated with the per-fund holdings table for that quarter.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: "总计" summary sheet - shift existing rows down and insert the
# new 2022-Q3 row at the top of the data (row 2).
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item(1)

# Push the existing data rows (2-8) down to (3-9), carrying values+styles.
$summary.Range("A2:D8").Copy($summary.Range("A3:D9"))

# Fill in the new first data row for 2022-Q3.
$summary.Range("A2").Value = 0
$summary.Range("B2").Value = "2022-Q3"
$summary.Range("C2").Value = 4
$summary.Range("D2").Value = 0.88

# Renumber the running index in column A (0..7) for all data rows.
$summary.Range("A3").Value = 1
$summary.Range("A4").Value = 2
$summary.Range("A5").Value = 3
$summary.Range("A6").Value = 4
$summary.Range("A7").Value = 5
$summary.Range("A8").Value = 6
$summary.Range("A9").Value = 7

# ---------------------------------------------------------------------
# Step 2: Create the new "2022-Q3" worksheet right after "总计" by
# duplicating the existing "2022-Q2" sheet (to inherit its formatting),
# then overwrite its contents with the 2022-Q3 fund holdings.
# ---------------------------------------------------------------------
$oldQ2 = $wb.Worksheets.Item(2)
# Copy(Before) inserts the duplicate immediately before $oldQ2 - i.e. right
# after "总计" - and pushes the original "2022-Q2" (and everything after
# it) one slot later.
$oldQ2.Copy($oldQ2)

$newSheet = $wb.Worksheets.Item(2)
$newSheet.Name = "2022-Q3"

# The source sheet had 5 data rows (rows 2-6); 2022-Q3 only has 4, so
# remove the extra row.
$newSheet.Rows.Item(6).Delete()

# Columns B and D:G hold values that look numeric ("001257", "68.35", ...)
# but must be stored as text, matching the rest of the workbook. Force
# text formatting before assigning so Excel doesn't coerce them to numbers.
$newSheet.Range("B2:B5").NumberFormat = "@"
$newSheet.Range("D2:G5").NumberFormat = "@"

$newSheet.Range("A2").Value = 0
$newSheet.Range("B2").Value = "001257"
$newSheet.Range("C2").Value = "兴业收益增强债券A"
$newSheet.Range("D2").Value = "68.35"
$newSheet.Range("E2").Value = "20.08"
$newSheet.Range("F2").Value = "0.74"
$newSheet.Range("G2").Value = "0.5058"
$newSheet.Range("H2").Value = 8

$newSheet.Range("A3").Value = 1
$newSheet.Range("B3").Value = "005984"
$newSheet.Range("C3").Value = "兴业聚华混合A"
$newSheet.Range("D3").Value = "23.61"
$newSheet.Range("E3").Value = "29.42"
$newSheet.Range("F3").Value = "0.96"
$newSheet.Range("G3").Value = "0.2267"
$newSheet.Range("H3").Value = 10

$newSheet.Range("A4").Value = 2
$newSheet.Range("B4").Value = "001258"
$newSheet.Range("C4").Value = "兴业收益增强债券C"
$newSheet.Range("D4").Value = "11.01"
$newSheet.Range("E4").Value = "20.08"
$newSheet.Range("F4").Value = "0.74"
$newSheet.Range("G4").Value = "0.0815"
$newSheet.Range("H4").Value = 8

$newSheet.Range("A5").Value = 3
$newSheet.Range("B5").Value = "005985"
$newSheet.Range("C5").Value = "兴业聚华混合C"
$newSheet.Range("D5").Value = "6.45"
$newSheet.Range("E5").Value = "29.42"
$newSheet.Range("F5").Value = "0.96"
$newSheet.Range("G5").Value = "0.0619"
$newSheet.Range("H5").Value = 10
